$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value looks numeric need to be forced to Text format first,
# otherwise Excel auto-converts the assigned string into a number (losing formatting
# like trailing zeros / the source's inline-string type). We set NumberFormat to "@"
# (Text), write the value, then restore the cell's style to "Normal" so no visible
# formatting change is left behind.

$ws.Range("D2").Value = "30.842.64"
$ws.Range("E2").Value = "  +1.47%  "
$ws.Range("D3").Value = "1.888.16"
$ws.Range("E3").Value = "  +2.16%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9981"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9985"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4777"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2883"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06582"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.91"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +15.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "97.93"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +16.73%  "
$ws.Range("D12").Value = "1.881.17"
$ws.Range("E12").Value = "  +1.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07578"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.138"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6578"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "313.04"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +36.66%  "
$ws.Range("D17").Value = "30.831.26"
$ws.Range("E17").Value = "  +1.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9991"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007603"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.00%  "
$ws.Range("D21").Value = "2.117.31"
$ws.Range("E21").Value = "  +2.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9991"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.136"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.183"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.328"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.98%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.40"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +13.90%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.955"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1071"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.361"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.180"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.985"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05048"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.36%  "
$ws.Range("E34").Value = "  +3.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7377"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.711"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01957"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.707"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.086"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9069"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "108.11"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.15%  "
$ws.Range("E42").Value = "  -0.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4226"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.81%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.655"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.96"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.380"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.35%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.083"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.09%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1227"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.84"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05620"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.42%  "
$ws.Range("E51").Value = "  +3.39%  "
